$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 201
$ws.Range("I18").Value = 201
$ws.Range("K18").Value = 201
$ws.Range("M18").Value = 83
$ws.Range("H19").Value = 798.7895
$ws.Range("I19").Value = 950
$ws.Range("J19").Value = 729
$ws.Range("K19").Value = 950
$ws.Range("L19").Value = 729
$ws.Range("M19").Value = -775
$ws.Range("N19").Value = -1079
$ws.Range("H21").Value = 10000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 10000
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -10936
$ws.Range("H23").Value = 10000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 10000
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -10468
$ws.Range("H26").Value = 10000
$ws.Range("J26").Value = 10000
$ws.Range("L26").Value = 10000
$ws.Range("N26").Value = -10688
$ws.Range("H32").Value = 567.3333
$ws.Range("I32").Value = 157.5
$ws.Range("J32").Value = 716.36365
$ws.Range("K32").Value = 157.5
$ws.Range("L32").Value = 716.36365
$ws.Range("M32").Value = 168.5
$ws.Range("N32").Value = -1368.36365
$ws.Range("H33").Value = 185.81482
$ws.Range("I33").Value = 179.85715
$ws.Range("J33").Value = 206.66667
$ws.Range("K33").Value = 179.85715
$ws.Range("L33").Value = 206.66667
$ws.Range("M33").Value = 49.14285000000001
$ws.Range("N33").Value = -664.6666700000001
$ws.Range("H39").Value = 601.2308
$ws.Range("I39").Value = 209.75
$ws.Range("J39").Value = 775.2222
$ws.Range("K39").Value = 629.25
$ws.Range("L39").Value = 2325.6666
$ws.Range("M39").Value = -333.25
$ws.Range("N39").Value = -2917.6666
$ws.Range("H51").Value = 100000.664
$ws.Range("I51").Value = 100001
$ws.Range("J51").Value = 100000
$ws.Range("K51").Value = 100001
$ws.Range("L51").Value = 100000
$ws.Range("M51").Value = -99517
$ws.Range("N51").Value = -100968
$ws.Range("H88").Value = 22819260
$ws.Range("J88").Value = 30424948
$ws.Range("L88").Value = 30424948
$ws.Range("N88").Value = -30425760
$ws.Range("H91").Value = 22819260
$ws.Range("J91").Value = 30424948
$ws.Range("L91").Value = 30424948
$ws.Range("N91").Value = -30427756
$ws.Range("H111").Value = 1100
$ws.Range("I111").Value = 966.6667
$ws.Range("K111").Value = 2900.0001
$ws.Range("M111").Value = 166.9998999999998
$ws.Range("H125").Value = 2133.2307
$ws.Range("I125").Value = 1207
$ws.Range("K125").Value = 10863
$ws.Range("M125").Value = -8403

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 24797
$ws.Range("J52").Value = 24797
$ws.Range("L52").Value = 24797
$ws.Range("N52").Value = -25433

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1403.5186
$ws.Range("I99").Value = 1135
$ws.Range("J99").Value = 1860
$ws.Range("K99").Value = 1135
$ws.Range("L99").Value = 1860
$ws.Range("M99").Value = 363
$ws.Range("N99").Value = -4856

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1266137.5
$ws.Range("I31").Value = 1736940.6
$ws.Range("J31").Value = 10662.5
$ws.Range("K31").Value = 1736940.6
$ws.Range("L31").Value = 10662.5
$ws.Range("M31").Value = -1736645.6
$ws.Range("N31").Value = -11252.5
$ws.Range("H34").Value = 1266137.5
$ws.Range("I34").Value = 1736940.6
$ws.Range("J34").Value = 10662.5
$ws.Range("K34").Value = 1736940.6
$ws.Range("L34").Value = 10662.5
$ws.Range("M34").Value = -1736738.6
$ws.Range("N34").Value = -11066.5
$ws.Range("H99").Value = 250003000
$ws.Range("I99").Value = 333336000
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 333336000
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -333334502
$ws.Range("N99").Value = -6996
$ws.Range("H126").Value = 250003000
$ws.Range("I126").Value = 333336000
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 1000008000
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -1000005530
$ws.Range("N126").Value = -16940

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 36462116
$ws.Range("I5").Value = 55555910
$ws.Range("J5").Value = 10327.272
$ws.Range("K5").Value = 166667730
$ws.Range("L5").Value = 30981.816
$ws.Range("M5").Value = -166667618
$ws.Range("N5").Value = -31205.816
$ws.Range("H121").Value = 704.4545000000001
$ws.Range("I121").Value = 362.25
$ws.Range("J121").Value = 900
$ws.Range("K121").Value = 1086.75
$ws.Range("L121").Value = 2700
$ws.Range("M121").Value = 223.25
$ws.Range("N121").Value = -5320
$ws.Range("H131").Value = 815.47
$ws.Range("I131").Value = 657.5
$ws.Range("J131").Value = 822.05206
$ws.Range("K131").Value = 1972.5
$ws.Range("L131").Value = 2466.15618
$ws.Range("M131").Value = 3067.5
$ws.Range("N131").Value = -12546.15618
$ws.Range("H135").Value = 36462116
$ws.Range("I135").Value = 55555910
$ws.Range("J135").Value = 10327.272
$ws.Range("K135").Value = 500003190
$ws.Range("L135").Value = 92945.448
$ws.Range("M135").Value = -500000655
$ws.Range("N135").Value = -98015.448

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3967.1667
$ws.Range("I7").Value = 3004
$ws.Range("J7").Value = 4159.8
$ws.Range("K7").Value = 3004
$ws.Range("L7").Value = 4159.8
$ws.Range("M7").Value = -2892
$ws.Range("N7").Value = -4383.8
$ws.Range("H40").Value = 35718100
$ws.Range("I40").Value = 4340
$ws.Range("K40").Value = 4340
$ws.Range("M40").Value = -4204
$ws.Range("H126").Value = 3967.1667
$ws.Range("I126").Value = 3004
$ws.Range("J126").Value = 4159.8
$ws.Range("K126").Value = 9012
$ws.Range("L126").Value = 12479.4
$ws.Range("M126").Value = -6542
$ws.Range("N126").Value = -17419.4

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1272.8334
$ws.Range("I81").Value = 1052.5714
$ws.Range("K81").Value = 2105.1428
$ws.Range("M81").Value = -1044.1428
$ws.Range("H84").Value = 1272.8334
$ws.Range("I84").Value = 1052.5714
$ws.Range("K84").Value = 10525.714
$ws.Range("M84").Value = -5221.714
